# fix bug about SalesRecord
# - Customer sheet: correct row 8 (id 7 "Daven" -> id 8 "Ivan") and add a
#   new row 9 for a second customer that shares the same id (8, "test").
# - Room sheet: row 3 (波斯湾) booking data corrected - price/maxCapacity
#   stored as text, bookUserid reset to 0, bookDate moved forward to
#   2022-12-11.
# - Active tab moves from SalesRecord to Room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Customer sheet
# ---------------------------------------------------------------------
$customer = $wb.Worksheets.Item("Customer")

# Row 8: was id=7 "Daven"/"SECRET"/13117826002 -> now id=8 "Ivan"/"null"/13188880000
$customer.Range("A8").Value = "'8"
$customer.Range("B8").Value = "Ivan"
$customer.Range("C8").Value = "null"
$customer.Range("G8").Value = "'13188880000"

# Reset style on the numeric-looking text cells so they don't pick up an
# accidental quote-prefix style (copy the plain/no-style format from a
# sibling cell in the same row that already carries no special style).
$customer.Range("B8").Copy() | Out-Null
$customer.Range("A8").PasteSpecial(-4122) | Out-Null
$customer.Range("G8").PasteSpecial(-4122) | Out-Null

# New row 9: second customer record re-using id 8 ("test" / MALE)
$customer.Range("A9").Value = "'8"
$customer.Range("B9").Value = "test"
$customer.Range("C9").Value = "MALE"
$customer.Range("D9").Value = "'0"
$customer.Range("E9").Value = "'0"
$customer.Range("F9").Value = "'123456"
$customer.Range("G9").Value = "'13177788945"

$customer.Range("B9").Copy() | Out-Null
$customer.Range("A9").PasteSpecial(-4122) | Out-Null
$customer.Range("D9").PasteSpecial(-4122) | Out-Null
$customer.Range("E9").PasteSpecial(-4122) | Out-Null
$customer.Range("F9").PasteSpecial(-4122) | Out-Null
$customer.Range("G9").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Room sheet
# ---------------------------------------------------------------------
$room = $wb.Worksheets.Item("Room")

$room.Range("A3").Value = "'2"
$room.Range("C3").Value = "'288.0"
$room.Range("D3").Value = 44906
$room.Range("E3").Value = 0
$room.Range("F3").Value = "'15"

$room.Range("B3").Copy() | Out-Null
$room.Range("A3").PasteSpecial(-4122) | Out-Null
$room.Range("C3").PasteSpecial(-4122) | Out-Null
$room.Range("F3").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Active tab moves from SalesRecord to Room
# ---------------------------------------------------------------------
$room.Activate() | Out-Null
$room.Range("E11").Select() | Out-Null
